$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The existing "_GoBack" bookmark (currently sitting right after the
#    second picture / before "输入cin>>") needs to move down to the end
#    of the document. Remove it from its current location first.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) The last paragraph in the document is empty. Append the copyright
#    notice to it as three runs (the middle one carries the
#    rFonts/lang formatting used for the author's name), then put the
#    "_GoBack" bookmark back around the very end of that text.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraRange = $lastPara.Range

# Target only the text inside the paragraph - not its trailing
# paragraph mark - so the insert lands *inside* the existing paragraph
# (keeping its pPr) instead of splitting it into a new one.
$insertRange = $d.Range($paraRange.Start, $paraRange.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Copyright &#169;2021-2099 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>HanxiaoZhang</w:t></w:r><w:r><w:t>. All rights reserved</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertRange.InsertXML($xml)
